$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K19").Value = 861.5
$ws.Range("M19").Value = -686.5
$ws.Range("H19").Value = 1626.8462
$ws.Range("I19").Value = 861.5
$ws.Range("K40").Value = 1071.125
$ws.Range("N40").Value = -1804.5454
$ws.Range("L40").Value = 1454.5454
$ws.Range("I40").Value = 1071.125
$ws.Range("H40").Value = 1227.3334
$ws.Range("M40").Value = -896.125
$ws.Range("J40").Value = 1454.5454
$ws.Range("I112").Value = 780
$ws.Range("M112").Value = -1232
$ws.Range("K112").Value = 2340
$ws.Range("H112").Value = 2403.697
$ws.Range("L112").Value = 8080.928400000001
$ws.Range("J112").Value = 2693.6428
$ws.Range("N112").Value = -10296.9284
$ws.Range("I116").Value = 3820
$ws.Range("H116").Value = 4327.8184
$ws.Range("N116").Value = -11635
$ws.Range("K116").Value = 3820
$ws.Range("J116").Value = 4751
$ws.Range("L116").Value = 4751
$ws.Range("M116").Value = -378
$ws.Range("H132").Value = 9263775
$ws.Range("I132").Value = 11907464
$ws.Range("K132").Value = 35722392
$ws.Range("M132").Value = -35719862
$ws.Range("L132").Value = 32587.875
$ws.Range("N132").Value = -37647.875
$ws.Range("J132").Value = 10862.625
$ws.Range("J137").Value = 1506.8462
$ws.Range("M137").Value = -678
$ws.Range("I137").Value = 1076
$ws.Range("H137").Value = 1236.0286
$ws.Range("L137").Value = 4520.5386
$ws.Range("N137").Value = -9620.5386
$ws.Range("K137").Value = 3228
$ws.Range("N138").Value = -16859.573
$ws.Range("L138").Value = 6579.572999999999
$ws.Range("J138").Value = 2193.191
$ws.Range("H138").Value = 2060.19

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N32").ClearContents()
$ws.Range("M32").Value = -2461.75
$ws.Range("K32").Value = 2748.75
$ws.Range("H32").Value = 2748.75
$ws.Range("I32").Value = 2748.75
$ws.Range("L32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("H35").Value = 1882
$ws.Range("M35").Value = -1476
$ws.Range("K35").Value = 1882
$ws.Range("I35").Value = 1882
$ws.Range("M36").ClearContents()
$ws.Range("I36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M45").Value = -1464.8462
$ws.Range("I45").Value = 1841.8462
$ws.Range("K45").Value = 1841.8462
$ws.Range("H45").Value = 2222.2
$ws.Range("N53").ClearContents()
$ws.Range("M53").ClearContents()
$ws.Range("L53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("I61").Value = 1500
$ws.Range("K61").Value = 1500
$ws.Range("H61").Value = 2207.6
$ws.Range("M61").Value = -1288
$ws.Range("N108").Value = -48680
$ws.Range("L108").Value = 41000
$ws.Range("H108").Value = 41000
$ws.Range("J108").Value = 41000
$ws.Range("M132").Value = -10658.666
$ws.Range("I132").Value = 4396.222
$ws.Range("H132").Value = 4004.875
$ws.Range("K132").Value = 13188.666
$ws.Range("M136").Value = -1950
$ws.Range("K136").Value = 4500
$ws.Range("I136").Value = 1500
$ws.Range("H136").Value = 2207.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J86").Value = 4252.5
$ws.Range("N86").Value = -6498.5
$ws.Range("K86").Value = 4431.615
$ws.Range("H86").Value = 4407.7334
$ws.Range("M86").Value = -3308.615
$ws.Range("L86").Value = 4252.5
$ws.Range("I86").Value = 4431.615
$ws.Range("H89").Value = 4407.7334
$ws.Range("J89").Value = 4252.5
$ws.Range("N89").Value = -32494.5
$ws.Range("I89").Value = 4431.615
$ws.Range("M89").Value = -16542.075
$ws.Range("K89").Value = 22158.075
$ws.Range("L89").Value = 21262.5
$ws.Range("M107").Value = 253
$ws.Range("J107").Value = 2971
$ws.Range("N107").Value = -6811
$ws.Range("K107").Value = 1667
$ws.Range("I107").Value = 1667
$ws.Range("H107").Value = 2101.6667
$ws.Range("L107").Value = 2971
$ws.Range("I134").Value = 9634.076999999999
$ws.Range("M134").Value = -26367.231
$ws.Range("H134").Value = 12062.1
$ws.Range("K134").Value = 28902.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7917.222
$ws.Range("I132").Value = 11259.8
$ws.Range("K132").Value = 33779.39999999999
$ws.Range("M132").Value = -31249.39999999999
$ws.Range("L132").Value = 11217
$ws.Range("N132").Value = -16277
$ws.Range("J132").Value = 3739
$ws.Range("M141").ClearContents()
$ws.Range("J141").Value = 33106
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 33106
$ws.Range("H141").Value = 33106
$ws.Range("I141").Value = 0
$ws.Range("N141").Value = -43466

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J114").Value = 1140.375
$ws.Range("L114").Value = 3421.125
$ws.Range("N114").Value = -9929.125
$ws.Range("H114").Value = 663.7
$ws.Range("K131").Value = 300001830
$ws.Range("I131").Value = 100000610
$ws.Range("N131").Value = -16139.7
$ws.Range("H131").Value = 33334882
$ws.Range("J131").Value = 2019.9
$ws.Range("M131").Value = -299996790
$ws.Range("L131").Value = 6059.700000000001
$ws.Range("J137").Value = 31008.25
$ws.Range("M137").Value = -7960.000499999998
$ws.Range("I137").Value = 4353.3335
$ws.Range("H137").Value = 19584.715
$ws.Range("L137").Value = 93024.75
$ws.Range("N137").Value = -103224.75
$ws.Range("K137").Value = 13060.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M126").Value = -2951.6
$ws.Range("K126").Value = 5421.6
$ws.Range("I126").Value = 1807.2
$ws.Range("H126").Value = 2721.4546
$ws.Range("M132").Value = -4468.0772
$ws.Range("I132").Value = 2332.6924
$ws.Range("H132").Value = 2673.4
$ws.Range("K132").Value = 6998.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 10000
$ws.Range("I87").Value = 10000
$ws.Range("M87").Value = -8877
$ws.Range("K87").Value = 10000
$ws.Range("M90").Value = -24384
$ws.Range("I90").Value = 10000
$ws.Range("H90").Value = 10000
$ws.Range("K90").Value = 30000
$ws.Range("M122").Value = -125004100
$ws.Range("J122").Value = 2430.8333
$ws.Range("I122").Value = 41668850
$ws.Range("N122").Value = -12192.4999
$ws.Range("K122").Value = 125006550
$ws.Range("H122").Value = 20835638
$ws.Range("L122").Value = 7292.499899999999
$ws.Range("M132").Value = -3193.3079
$ws.Range("I132").Value = 1907.7693
$ws.Range("H132").Value = 52310.45
$ws.Range("K132").Value = 5723.3079
$ws.Range("M136").Value = -4402.200000000001
$ws.Range("I136").Value = 2317.4
$ws.Range("J136").Value = 1558.4286
$ws.Range("K136").Value = 6952.200000000001
$ws.Range("H136").Value = 2004.8823
$ws.Range("L136").Value = 4675.2858
$ws.Range("N136").Value = -9775.2858
$ws.Range("H137").Value = 49143
$ws.Range("L137").Value = 49143
$ws.Range("N137").Value = -59343
$ws.Range("J137").Value = 49143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I96").Value = 1642.75
$ws.Range("K96").Value = 1642.75
$ws.Range("M96").Value = -269.75
$ws.Range("H96").Value = 1812.2
$ws.Range("M132").Value = -3105.799999999999
$ws.Range("I132").Value = 1878.6
$ws.Range("H132").Value = 2260.1924
$ws.Range("K132").Value = 5635.799999999999
$ws.Range("M136").Value = -2025.75
$ws.Range("H136").Value = 1767.0834
$ws.Range("I136").Value = 1525.25
$ws.Range("K136").Value = 4575.75
